$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 81.05837566666666
$ws.Range("H2").Value = 243.175127
$ws.Range("I2").Value = 0.3545816884225585
$ws.Range("J2").Value = 0.3545816884225585
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 8020.460724977983
$ws.Range("R2").Value = 72184.14652480184
$ws.Range("S2").Value = 0.07439158192672862
$ws.Range("T2").Value = 0.07439158192672864

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 81.05837566666666
$ws.Range("H3").Value = 243.175127
$ws.Range("I3").Value = 0.3545816884225585
$ws.Range("J3").Value = 0.3545816884225585
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 13213.02068667786
$ws.Range("R3").Value = 118917.1861801008
$ws.Range("S3").Value = 0.122553746551169
$ws.Range("T3").Value = 0.122553746551169

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 81.05837566666666
$ws.Range("H4").Value = 243.175127
$ws.Range("I4").Value = 0.3545816884225585
$ws.Range("J4").Value = 0.3545816884225585
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 5300.902181324072
$ws.Range("R4").Value = 47708.11963191664
$ws.Range("S4").Value = 0.04916706314382292
$ws.Range("T4").Value = 0.04916706314382293

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 81.05837566666666
$ws.Range("H5").Value = 243.175127
$ws.Range("I5").Value = 0.3545816884225585
$ws.Range("J5").Value = 0.3545816884225585
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 11694.51855068728
$ws.Range("R5").Value = 105250.6669561855
$ws.Range("S5").Value = 0.1084692968008379
$ws.Range("T5").Value = 0.108469296800838

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 123.018252
$ws.Range("H6").Value = 369.054756
$ws.Range("I6").Value = 0.5381309351710768
$ws.Range("J6").Value = 0.5381309351710768
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 12172.25302760645
$ws.Range("R6").Value = 109550.277248458
$ws.Range("S6").Value = 0.112900391808669
$ws.Range("T6").Value = 0.112900391808669

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 123.018252
$ws.Range("H7").Value = 369.054756
$ws.Range("I7").Value = 0.5381309351710768
$ws.Range("J7").Value = 0.5381309351710768
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 20052.74217681339
$ws.Range("R7").Value = 180474.6795913205
$ws.Range("S7").Value = 0.1859937058049835
$ws.Range("T7").Value = 0.1859937058049836

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 123.018252
$ws.Range("H8").Value = 369.054756
$ws.Range("I8").Value = 0.5381309351710768
$ws.Range("J8").Value = 0.5381309351710768
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 8044.914729738882
$ws.Range("R8").Value = 72404.23256764993
$ws.Range("S8").Value = 0.07461839833553438
$ws.Range("T8").Value = 0.07461839833553439

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 123.018252
$ws.Range("H9").Value = 369.054756
$ws.Range("I9").Value = 0.5381309351710768
$ws.Range("J9").Value = 0.5381309351710768
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 17748.18725708476
$ws.Range("R9").Value = 159733.6853137628
$ws.Range("S9").Value = 0.1646184392218898
$ws.Range("T9").Value = 0.1646184392218898

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3624666666666667
$ws.Range("H10").Value = 1.0874
$ws.Range("I10").Value = 0.001585573873230423
$ws.Range("J10").Value = 0.001585573873230423
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 35.86488922586667
$ws.Range("R10").Value = 322.7840030327999
$ws.Range("S10").Value = 0.0003326549354989119
$ws.Range("T10").Value = 0.0003326549354989119

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3624666666666667
$ws.Range("H11").Value = 1.0874
$ws.Range("I11").Value = 0.001585573873230423
$ws.Range("J11").Value = 0.001585573873230423
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 59.08432688797777
$ws.Range("R11").Value = 531.7589419917999
$ws.Range("S11").Value = 0.0005480204560548708
$ws.Range("T11").Value = 0.0005480204560548709

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3624666666666667
$ws.Range("H12").Value = 1.0874
$ws.Range("I12").Value = 0.001585573873230423
$ws.Range("J12").Value = 0.001585573873230423
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 23.70390879644445
$ws.Range("R12").Value = 213.335179168
$ws.Range("S12").Value = 0.0002198590995804972
$ws.Range("T12").Value = 0.0002198590995804972

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3624666666666667
$ws.Range("H13").Value = 1.0874
$ws.Range("I13").Value = 0.001585573873230423
$ws.Range("J13").Value = 0.001585573873230423
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 52.29407969844444
$ws.Range("R13").Value = 470.646717286
$ws.Range("S13").Value = 0.0004850393820961433
$ws.Range("T13").Value = 0.0004850393820961434

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 24.16373066666667
$ws.Range("H14").Value = 72.491192
$ws.Range("I14").Value = 0.1057018025331343
$ws.Range("J14").Value = 0.1057018025331344
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 2390.921989085003
$ws.Range("R14").Value = 21518.29790176502
$ws.Range("S14").Value = 0.02217634062810303
$ws.Range("T14").Value = 0.02217634062810304

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 24.16373066666667
$ws.Range("H15").Value = 72.491192
$ws.Range("I15").Value = 0.1057018025331343
$ws.Range("J15").Value = 0.1057018025331344
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 3938.838775636527
$ws.Range("R15").Value = 35449.54898072874
$ws.Range("S15").Value = 0.03653361789571566
$ws.Range("T15").Value = 0.03653361789571567

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 24.16373066666667
$ws.Range("H16").Value = 72.491192
$ws.Range("I16").Value = 0.1057018025331343
$ws.Range("J16").Value = 0.1057018025331344
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 1580.213908141938
$ws.Range("R16").Value = 14221.92517327744
$ws.Range("S16").Value = 0.01465684035372167
$ws.Range("T16").Value = 0.01465684035372167

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 24.16373066666667
$ws.Range("H17").Value = 72.491192
$ws.Range("I17").Value = 0.1057018025331343
$ws.Range("J17").Value = 0.1057018025331344
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 3486.169001180097
$ws.Range("R17").Value = 31375.52101062088
$ws.Range("S17").Value = 0.03233500365559398
$ws.Range("T17").Value = 0.03233500365559398
